# Word COM-interop script implementing the "add words in June 7th" commit.
#
# Changes applied:
#  1. The existing paragraph's rFonts hint is normalized from
#     hint="default" eastAsiaTheme="minorEastAsia" to hint="eastAsia".
#  2. A new paragraph "今天天气不错！心情也很好！" is appended after the
#     existing one, and the "_GoBack" bookmark (which Word re-stamps on the
#     last edit point) moves to the end of that new paragraph.
#  3. Two built-in styles ("Default Paragraph Font" / "Normal Table") are
#     flipped to Quick Style (w:qFormat) to match the authoring tool's
#     re-save of styles.xml.

$d = $word.ActiveDocument

# --- 1 & 2: rebuild the document body with the exact run/paragraph
#            properties the author's Word session produced. Using
#            InsertXML on the whole body lets us set w:rFonts/w:hint and
#            the bookmark position precisely, instead of relying on
#            whatever defaults Word's higher-level paragraph APIs pick.
$body = $d.Content

$openXml = '<?xml version="1.0"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
'<pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
'<w:body>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:rPr>' + `
      '<w:rFonts w:hint="eastAsia"/>' + `
      '<w:lang w:val="en-US" w:eastAsia="zh-CN"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:hint="eastAsia"/>' + `
      '<w:lang w:val="en-US" w:eastAsia="zh-CN"/>' + `
    '</w:rPr>' + `
    '<w:t>蛐蛐很难受！</w:t>' + `
  '</w:r>' + `
'</w:p>' + `
'<w:p>' + `
  '<w:pPr>' + `
    '<w:rPr>' + `
      '<w:rFonts w:hint="default"/>' + `
      '<w:lang w:val="en-US" w:eastAsia="zh-CN"/>' + `
    '</w:rPr>' + `
  '</w:pPr>' + `
  '<w:r>' + `
    '<w:rPr>' + `
      '<w:rFonts w:hint="eastAsia"/>' + `
      '<w:lang w:val="en-US" w:eastAsia="zh-CN"/>' + `
    '</w:rPr>' + `
    '<w:t>今天天气不错！心情也很好！</w:t>' + `
  '</w:r>' + `
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
  '<w:bookmarkEnd w:id="0"/>' + `
'</w:p>' + `
'</w:body>' + `
'</w:document>' + `
'</pkg:xmlData>' + `
'</pkg:part>' + `
'</pkg:package>'

[void]$body.InsertXML($openXml)

# --- 3: mark the two latent "Default Paragraph Font" / "Normal Table"
#        styles as Quick Styles (adds <w:qFormat/> to their definitions).
$d.Styles("Default Paragraph Font").QuickStyle = $true
$d.Styles("Normal Table").QuickStyle = $true

Write-Output "edit applied"
